$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their exact literal representation (avoid Excel
# auto-converting numeric-looking strings into floating point numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.387.39'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -9.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.866.15'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -9.28%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '544.88'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -8.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '120.05'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -11.34%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.862.26'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -9.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.489'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -4.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.124'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -12.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '4.77'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -11.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.426'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -5.72%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000210'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -11.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '31.08'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -10.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.118'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.340.09'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -9.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.873.62'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -9.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '57.476.77'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -9.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.25'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -4.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '410.73'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -10.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.67'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -8.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.645'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -7.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.70'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -12.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.40'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -6.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '76.07'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -8.32%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.42'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -10.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.88'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -9.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.01'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -8.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '24.40'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -9.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.86'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -12.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0922'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -7.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.31'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -9.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '48.19'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.881'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -13.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.95'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -18.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.27'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0607'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -16.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0338'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -12.96%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -7.68%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.592.28'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -6.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '345.98'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -11.15%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.29'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -12.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '117.44'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -6.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.225'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -9.78%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.90'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -9.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.42'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -10.07%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -11.17%  '
